$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 519, shifting existing rows 519:544 down to 520:545
$ws.Rows(519).Insert()

# Populate the new row 519 with the new record's data
$ws.Cells.Item(519, 1).Value = 5
$ws.Cells.Item(519, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(519, 3).Value = "Maule"
$ws.Cells.Item(519, 4).Value = 45041
$ws.Cells.Item(519, 5).Value = 7
$ws.Cells.Item(519, 6).Value = 100112032
$ws.Cells.Item(519, 7).Value = "Zapallo italiano"
$ws.Cells.Item(519, 8).Value = "Sin especificar"
$ws.Cells.Item(519, 9).Value = "Primera"
$ws.Cells.Item(519, 10).Value = 200
$ws.Cells.Item(519, 11).Value = 11000
$ws.Cells.Item(519, 12).Value = 11000
$ws.Cells.Item(519, 13).Value = 11000
$ws.Cells.Item(519, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(519, 15).Value = "Región del Maule"
$ws.Cells.Item(519, 16).Value = 220
$ws.Cells.Item(519, 17).Value = 50
$ws.Cells.Item(519, 18).Value = "Hortaliza"
